$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.36022366666667
$ws.Range("H2").Value = 58.080671
$ws.Range("I2").Value = 0.005884129141485179
$ws.Range("J2").Value = 0.005884129141485179
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 169.629438
$ws.Range("N2").Value = 508.888314
$ws.Range("O2").Value = 0.7428377317484701
$ws.Range("P2").Value = 0.7428377317484702
$ws.Range("Q2").Value = 3284.063860130966
$ws.Range("R2").Value = 29556.57474117869
$ws.Range("S2").Value = 0.004370953144775923
$ws.Range("T2").Value = 0.004370953144775923

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.36022366666667
$ws.Range("H3").Value = 58.080671
$ws.Range("I3").Value = 0.005884129141485179
$ws.Range("J3").Value = 0.005884129141485179
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.9848756666666668
$ws.Range("N3").Value = 2.954627
$ws.Range("O3").Value = 0.004312947180081616
$ws.Range("P3").Value = 0.004312947180081616
$ws.Range("Q3").Value = 19.06741319052411
$ws.Range("R3").Value = 171.606718714717
$ws.Range("S3").Value = 0.00002537793818800456
$ws.Range("T3").Value = 0.00002537793818800456

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.36022366666667
$ws.Range("H4").Value = 58.080671
$ws.Range("I4").Value = 0.005884129141485179
$ws.Range("J4").Value = 0.005884129141485179
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.620752
$ws.Range("N4").Value = 163.862256
$ws.Range("O4").Value = 0.2391940691454494
$ws.Range("P4").Value = 0.2391940691454494
$ws.Range("Q4").Value = 1057.469975561531
$ws.Range("R4").Value = 9517.229780053774
$ws.Range("S4").Value = 0.00140744879272916
$ws.Range("T4").Value = 0.00140744879272916

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.36022366666667
$ws.Range("H5").Value = 58.080671
$ws.Range("I5").Value = 0.005884129141485179
$ws.Range("J5").Value = 0.005884129141485179
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.118221666666667
$ws.Range("N5").Value = 9.354665000000001
$ws.Range("O5").Value = 0.01365525192599884
$ws.Range("P5").Value = 0.01365525192599884
$ws.Range("Q5").Value = 60.36946890891278
$ws.Range("R5").Value = 543.325220180215
$ws.Range("S5").Value = 0.00008034926579209139
$ws.Range("T5").Value = 0.00008034926579209139

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3161.845459
$ws.Range("H6").Value = 9485.536377
$ws.Range("I6").Value = 0.9609758299542277
$ws.Range("J6").Value = 0.9609758299542278
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 169.629438
$ws.Range("N6").Value = 508.888314
$ws.Range("O6").Value = 0.7428377317484701
$ws.Range("P6").Value = 0.7428377317484702
$ws.Range("Q6").Value = 536342.068253022
$ws.Range("R6").Value = 4827078.614277198
$ws.Range("S6").Value = 0.7138491057883021
$ws.Range("T6").Value = 0.7138491057883022

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3161.845459
$ws.Range("H7").Value = 9485.536377
$ws.Range("I7").Value = 0.9609758299542277
$ws.Range("J7").Value = 0.9609758299542278
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9848756666666668
$ws.Range("N7").Value = 2.954627
$ws.Range("O7").Value = 0.004312947180081616
$ws.Range("P7").Value = 0.004312947180081616
$ws.Range("Q7").Value = 3114.024654329598
$ws.Range("R7").Value = 28026.22188896638
$ws.Range("S7").Value = 0.004144637995927677
$ws.Range("T7").Value = 0.004144637995927677

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3161.845459
$ws.Range("H8").Value = 9485.536377
$ws.Range("I8").Value = 0.9609758299542277
$ws.Range("J8").Value = 0.9609758299542278
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.620752
$ws.Range("N8").Value = 163.862256
$ws.Range("O8").Value = 0.2391940691454494
$ws.Range("P8").Value = 0.2391940691454494
$ws.Range("Q8").Value = 172702.3766783652
$ws.Range("R8").Value = 1554321.390105287
$ws.Range("S8").Value = 0.2298597191171771
$ws.Range("T8").Value = 0.2298597191171772

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3161.845459
$ws.Range("H9").Value = 9485.536377
$ws.Range("I9").Value = 0.9609758299542277
$ws.Range("J9").Value = 0.9609758299542278
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.118221666666667
$ws.Range("N9").Value = 9.354665000000001
$ws.Range("O9").Value = 0.01365525192599884
$ws.Range("P9").Value = 0.01365525192599884
$ws.Range("Q9").Value = 9859.335016905412
$ws.Range("R9").Value = 88734.01515214871
$ws.Range("S9").Value = 0.0131223670528208
$ws.Range("T9").Value = 0.0131223670528208

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.055785333333333
$ws.Range("H10").Value = 6.167356
$ws.Range("I10").Value = 0.0006248123263850286
$ws.Range("J10").Value = 0.0006248123263850286
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 169.629438
$ws.Range("N10").Value = 508.888314
$ws.Range("O10").Value = 0.7428377317484701
$ws.Range("P10").Value = 0.7428377317484702
$ws.Range("Q10").Value = 348.7217107419759
$ws.Range("R10").Value = 3138.495396677784
$ws.Range("S10").Value = 0.0004641341713003394
$ws.Range("T10").Value = 0.0004641341713003395

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.055785333333333
$ws.Range("H11").Value = 6.167356
$ws.Range("I11").Value = 0.0006248123263850286
$ws.Range("J11").Value = 0.0006248123263850286
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.9848756666666668
$ws.Range("N11").Value = 2.954627
$ws.Range("O11").Value = 0.004312947180081616
$ws.Range("P11").Value = 0.004312947180081616
$ws.Range("Q11").Value = 2.024692950690222
$ws.Range("R11").Value = 18.222236556212
$ws.Range("S11").Value = 0.000002694782561162543
$ws.Range("T11").Value = 0.000002694782561162543

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.055785333333333
$ws.Range("H12").Value = 6.167356
$ws.Range("I12").Value = 0.0006248123263850286
$ws.Range("J12").Value = 0.0006248123263850286
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.620752
$ws.Range("N12").Value = 163.862256
$ws.Range("O12").Value = 0.2391940691454494
$ws.Range("P12").Value = 0.2391940691454494
$ws.Range("Q12").Value = 112.2885408572373
$ws.Range("R12").Value = 1010.596867715136
$ws.Range("S12").Value = 0.0001494514028002696
$ws.Range("T12").Value = 0.0001494514028002696

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.055785333333333
$ws.Range("H13").Value = 6.167356
$ws.Range("I13").Value = 0.0006248123263850286
$ws.Range("J13").Value = 0.0006248123263850286
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.118221666666667
$ws.Range("N13").Value = 9.354665000000001
$ws.Range("O13").Value = 0.01365525192599884
$ws.Range("P13").Value = 0.01365525192599884
$ws.Range("Q13").Value = 6.410394368415555
$ws.Range("R13").Value = 57.69354931574
$ws.Range("S13").Value = 0.000008531969723256979
$ws.Range("T13").Value = 0.000008531969723256979

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 106.9830526666667
$ws.Range("H14").Value = 320.949158
$ws.Range("I14").Value = 0.03251522857790212
$ws.Range("J14").Value = 0.03251522857790212
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 169.629438
$ws.Range("N14").Value = 508.888314
$ws.Range("O14").Value = 0.7428377317484701
$ws.Range("P14").Value = 0.7428377317484702
$ws.Range("Q14").Value = 18147.47509937107
$ws.Range("R14").Value = 163327.2758943396
$ws.Range("S14").Value = 0.02415353864409184
$ws.Range("T14").Value = 0.02415353864409185

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 106.9830526666667
$ws.Range("H15").Value = 320.949158
$ws.Range("I15").Value = 0.03251522857790212
$ws.Range("J15").Value = 0.03251522857790212
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.9848756666666668
$ws.Range("N15").Value = 2.954627
$ws.Range("O15").Value = 0.004312947180081616
$ws.Range("P15").Value = 0.004312947180081616
$ws.Range("Q15").Value = 105.3650053171185
$ws.Range("R15").Value = 948.2850478540662
$ws.Range("S15").Value = 0.0001402364634047721
$ws.Range("T15").Value = 0.0001402364634047721

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 106.9830526666667
$ws.Range("H16").Value = 320.949158
$ws.Range("I16").Value = 0.03251522857790212
$ws.Range("J16").Value = 0.03251522857790212
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 54.620752
$ws.Range("N16").Value = 163.862256
$ws.Range("O16").Value = 0.2391940691454494
$ws.Range("P16").Value = 0.2391940691454494
$ws.Range("Q16").Value = 5843.494787908939
$ws.Range("R16").Value = 52591.45309118045
$ws.Range("S16").Value = 0.007777449832742812
$ws.Range("T16").Value = 0.007777449832742812

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 106.9830526666667
$ws.Range("H17").Value = 320.949158
$ws.Range("I17").Value = 0.03251522857790212
$ws.Range("J17").Value = 0.03251522857790212
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.118221666666667
$ws.Range("N17").Value = 9.354665000000001
$ws.Range("O17").Value = 0.01365525192599884
$ws.Range("P17").Value = 0.01365525192599884
$ws.Range("Q17").Value = 333.5968727913411
$ws.Range("R17").Value = 3002.37185512207
$ws.Range("S17").Value = 0.0004440036376626906
$ws.Range("T17").Value = 0.0004440036376626906
